$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 152.09
$ws.Range("C2").Value = 84.23
$ws.Range("B3").Value = 136.62
$ws.Range("C3").Value = 82.55
$ws.Range("B8").Value = 30

$ws.Range("F9:F15").Select()
